# Realizacao de um ecra de teste para o aluno. O teste sera efectuado
# atraves de uma associacao de um teste a uma imagem.
#
# Assign two pending tasks (on the "Folha1" task sheet) to students and
# move their status from "Por fazer" to "A fazer...".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

# Use an existing "A fazer.." status cell as the format template (orange fill)
$formatSource = $ws.Range("F9")

# Row 11 -> "Fazer o pedido do teste e recebe-lo (online)" assigned to Jorge 13683
$ws.Range("E11").Value = "Jorge 13683"
$ws.Range("F11").Value = "A fazer..."
$formatSource.Copy()
$ws.Range("F11").PasteSpecial(-4122)

# Row 15 -> "Executar o teste com imagens (offline) modo aluno" assigned to Rafael 13696
$ws.Range("E15").Value = "Rafael 13696"
$ws.Range("F15").Value = "A fazer..."
$formatSource.Copy()
$ws.Range("F15").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Update the active selection to reflect where the user ended up working
$ws.Range("H14").Select()
